$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Oct 09 22:48:23 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 22:48:37 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 22:48:51 EDT 2023"
